$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("itemloc")
$ws.Activate()

# The sheet has columns D (max capacity), E (expiry date), F (fifo date) over
# rows 1:63 (header + 62 data rows). The edit rotates these three columns one
# step to the left: new D = old E, new E = old F, new F = old D.
# Use Copy (not Value assignment) so literal text that looks like a date
# (e.g. "2100-01-01") is relocated as-is instead of being re-interpreted/
# coerced into a real date serial value, and use a scratch column instead of
# whole-column Cut/Insert so column widths aren't recalculated.

$scratch = $ws.Range("K1:K63")

$ws.Range("D1:D63").Copy($scratch)
$ws.Range("E1:E63").Copy($ws.Range("D1:D63"))
$ws.Range("F1:F63").Copy($ws.Range("E1:E63"))
$scratch.Copy($ws.Range("F1:F63"))
$scratch.Clear()

# Leave the UI selection on the whole of column F (where the old column D
# data now lives), matching the end state of the recorded session.
$ws.Columns.Item(6).Select()
